# Applies the updated Brynhildr_Profits Leve-profit figures (H/I/J/K/L/M/N columns)
# for the rows that moved in this scheduled-runner refresh, one worksheet at a time.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 3484.2
$ws.Range("J40").Value = 3815.3333
$ws.Range("L40").Value = 3815.3333
$ws.Range("N40").Value = -4165.3333

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 12680.857
$ws.Range("I132").Value = 12938.941
$ws.Range("K132").Value = 38816.823
$ws.Range("M132").Value = -36286.823

# Row 138: All-night Crafting
$ws.Range("H138").Value = 11099.8
$ws.Range("I138").Value = 11625
$ws.Range("J138").Value = 8999
$ws.Range("K138").Value = 34875
$ws.Range("L138").Value = 26997
$ws.Range("M138").Value = -29735
$ws.Range("N138").Value = -37277

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 1204.5
$ws.Range("I2").Value = 1087.8235
$ws.Range("K2").Value = 1087.8235
$ws.Range("M2").Value = -974.8235

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 1038920.75
$ws.Range("I32").Value = 1112208.9
$ws.Range("J32").Value = 37316.668
$ws.Range("K32").Value = 1112208.9
$ws.Range("L32").Value = 37316.668
$ws.Range("M32").Value = -1111921.9
$ws.Range("N32").Value = -37890.668

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3227740
$ws.Range("I61").Value = 2042.3334
$ws.Range("K61").Value = 2042.3334
$ws.Range("M61").Value = -1830.3334

# Row 116: No Scope
$ws.Range("H116").Value = 1204.5
$ws.Range("I116").Value = 1087.8235
$ws.Range("K116").Value = 1087.8235
$ws.Range("M116").Value = 1206.1765

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 4430.567
$ws.Range("I132").Value = 2272.0588
$ws.Range("J132").Value = 7253.231
$ws.Range("K132").Value = 6816.176399999999
$ws.Range("L132").Value = 21759.693
$ws.Range("M132").Value = -4286.176399999999
$ws.Range("N132").Value = -26819.693

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 3227740
$ws.Range("I136").Value = 2042.3334
$ws.Range("K136").Value = 6127.0002
$ws.Range("M136").Value = -3577.0002

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 1204.5
$ws.Range("I3").Value = 1087.8235
$ws.Range("K3").Value = 1087.8235
$ws.Range("M3").Value = -973.8235

# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 28918.105
$ws.Range("I20").Value = 33465.188
$ws.Range("J20").Value = 12093.9
$ws.Range("K20").Value = 33465.188
$ws.Range("L20").Value = 12093.9
$ws.Range("M20").Value = -33218.188
$ws.Range("N20").Value = -12587.9

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 2862.9333
$ws.Range("I86").Value = 4201
$ws.Range("J86").Value = 1970.8889
$ws.Range("K86").Value = 4201
$ws.Range("L86").Value = 1970.8889
$ws.Range("M86").Value = -3078
$ws.Range("N86").Value = -4216.8889

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 2862.9333
$ws.Range("I89").Value = 4201
$ws.Range("J89").Value = 1970.8889
$ws.Range("K89").Value = 21005
$ws.Range("L89").Value = 9854.4445
$ws.Range("M89").Value = -15389
$ws.Range("N89").Value = -21086.4445

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 3262.8125
$ws.Range("I105").Value = 1171.8572
$ws.Range("K105").Value = 1171.8572
$ws.Range("M105").Value = 575.1428000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 147042.72
$ws.Range("I16").Value = 4600.5
$ws.Range("K16").Value = 4600.5
$ws.Range("M16").Value = -4313.5

# Row 31: Wall Not Found
$ws.Range("H31").Value = 1086102.9
$ws.Range("I31").Value = 1291751.1
$ws.Range("K31").Value = 1291751.1
$ws.Range("M31").Value = -1291456.1

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 1086102.9
$ws.Range("I34").Value = 1291751.1
$ws.Range("K34").Value = 1291751.1
$ws.Range("M34").Value = -1291549.1

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 4618
$ws.Range("I62").Value = 4363.3335
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4363.3335
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3739.3335
$ws.Range("N62").Value = -6248

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 4618
$ws.Range("I65").Value = 4363.3335
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 21816.6675
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -18696.6675
$ws.Range("N65").Value = -31240

# Row 99: O Pine
$ws.Range("H99").Value = 1640799.6
$ws.Range("I99").Value = 54499.75
$ws.Range("J99").Value = 7985999
$ws.Range("K99").Value = 54499.75
$ws.Range("L99").Value = 7985999
$ws.Range("M99").Value = -53001.75
$ws.Range("N99").Value = -7988995

# Row 113: Patient Patients
$ws.Range("H113").Value = 147042.72
$ws.Range("I113").Value = 4600.5
$ws.Range("K113").Value = 4600.5
$ws.Range("M113").Value = -2430.5

# Row 126: A Better Conductor
$ws.Range("H126").Value = 1640799.6
$ws.Range("I126").Value = 54499.75
$ws.Range("J126").Value = 7985999
$ws.Range("K126").Value = 163499.25
$ws.Range("L126").Value = 23957997
$ws.Range("M126").Value = -161029.25
$ws.Range("N126").Value = -23962937

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1983.0851
$ws.Range("I132").Value = 1845.25
$ws.Range("K132").Value = 5535.75
$ws.Range("M132").Value = -3005.75

# Row 133: Yimepi's Country Charms
$ws.Range("H133").Value = 49880
$ws.Range("J133").Value = 49880
$ws.Range("L133").Value = 49880
$ws.Range("N133").Value = -54940

$ws = $wb.Worksheets.Item("CUL")
# Row 18: Fisher of Men
$ws.Range("H18").Value = 167760.5
$ws.Range("I18").Value = 1000000
$ws.Range("K18").Value = 3000000
$ws.Range("M18").Value = -2999831

# Row 139: Najoothie
$ws.Range("H139").Value = 5193.591
$ws.Range("I139").Value = 2891.1875
$ws.Range("J139").Value = 11333.333
$ws.Range("K139").Value = 8673.5625
$ws.Range("L139").Value = 33999.999
$ws.Range("M139").Value = -3533.5625
$ws.Range("N139").Value = -44279.999

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar
$ws.Range("H132").Value = 15909.75
$ws.Range("I132").Value = 9102.5
$ws.Range("K132").Value = 27307.5
$ws.Range("M132").Value = -24777.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 3144.3157
$ws.Range("I22").Value = 2727.3333
$ws.Range("J22").Value = 3519.6
$ws.Range("K22").Value = 2727.3333
$ws.Range("L22").Value = 3519.6
$ws.Range("M22").Value = -2432.3333
$ws.Range("N22").Value = -4109.6

# Row 27: Fire and Hide
$ws.Range("H27").Value = 3144.3157
$ws.Range("I27").Value = 2727.3333
$ws.Range("J27").Value = 3519.6
$ws.Range("K27").Value = 2727.3333
$ws.Range("L27").Value = 3519.6
$ws.Range("M27").Value = -2620.3333
$ws.Range("N27").Value = -3733.6

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 3195.3809
$ws.Range("J68").Value = 7801
$ws.Range("L68").Value = 7801
$ws.Range("N68").Value = -9299

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 3195.3809
$ws.Range("J71").Value = 7801
$ws.Range("L71").Value = 39005
$ws.Range("N71").Value = -46493

# Row 122: Hell on Leather
$ws.Range("H122").Value = 3458.3684
$ws.Range("I122").Value = 2983.1428
$ws.Range("J122").Value = 4789
$ws.Range("K122").Value = 8949.428400000001
$ws.Range("L122").Value = 14367
$ws.Range("M122").Value = -6499.428400000001
$ws.Range("N122").Value = -19267

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax
$ws.Range("H107").Value = 2151.2222
$ws.Range("I107").Value = 733.0625
$ws.Range("K107").Value = 2199.1875
$ws.Range("M107").Value = -279.1875

# Row 113: A Tender Table
$ws.Range("H113").Value = 1180
$ws.Range("J113").Value = 1683.1052
$ws.Range("L113").Value = 5049.3156
$ws.Range("N113").Value = -9389.3156
